$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before current row 2 (shifts everything, incl. formatting
# and formulas, down by 6 rows).
$ws.Rows("2:7").Insert()

# The newly inserted rows inherit formatting from row 8 (the old row 2) only
# for cells that had it; make column A match the rest of the date column by
# copying its number format down.
$ws.Cells.Item(8, 1).Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill the new rows: dates 2023-12-19 .. 2023-12-24 in column A, and the new
# "S" status flag in column C.
$dates = @(45279, 45280, 45281, 45282, 45283, 45284)
for ($i = 0; $i -lt 6; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 3).Value = "S"
}

# Update selection to match the committed state.
$ws.Range("A2:A8").Select()
$ws.Application.ActiveCell = $ws.Range("A8")

Write-Output "done"
